# Översikt MÖLNDAL.xlsx - automatic update of files
#
# Changes applied:
#  1. Column C ("Förändrad") bumped from 46059 -> 46060 for every data row (2..47).
#  2. The scraped dataset was refreshed, which re-shuffled several records to
#     different rows while keeping their own data intact:
#       - Row 5 and Row 6 fully swap their record content (including the
#         per-record HYPERLINK formulas in S/T/V/W/X/Y/Z).
#       - Rows 17..45 keep the same "shell" (Län/Kommun/species-count columns,
#         which are identical across these rows) but the Beteckning (A), Datum
#         (B) and Area (G) values move between rows according to the mapping
#         below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: bump "Förändrad" (column C) on every row from 46059 to 46060.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 47; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46059) {
        $cell.Value = 46060
    }
}

# ---------------------------------------------------------------------------
# Step 2: rows 17..45 - permute Beteckning/Datum/Area (columns A, B, G).
# Mapping is "new row" -> "old row" (i.e. new row N receives the A/B/G values
# that currently live on old row M).
# ---------------------------------------------------------------------------
$rowMap = [ordered]@{
    17 = 24
    18 = 20
    19 = 36
    20 = 22
    21 = 35
    22 = 39
    23 = 25
    24 = 18
    25 = 17
    26 = 19
    27 = 31
    28 = 32
    29 = 45
    30 = 44
    31 = 40
    32 = 42
    33 = 21
    34 = 26
    35 = 34
    36 = 30
    37 = 27
    38 = 29
    39 = 28
    40 = 33
    41 = 23
    42 = 43
    43 = 37
    44 = 41
    45 = 38
}

# Snapshot the current A/B/G values for every row involved before any of them
# get overwritten (several rows are both a source and a destination).
$snapshotA = @{}
$snapshotB = @{}
$snapshotG = @{}
foreach ($oldRow in $rowMap.Values) {
    $snapshotA[$oldRow] = $ws.Cells.Item($oldRow, 1).Value2
    $snapshotB[$oldRow] = $ws.Cells.Item($oldRow, 2).Value2
    $snapshotG[$oldRow] = $ws.Cells.Item($oldRow, 7).Value2
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $ws.Cells.Item($newRow, 1).Value = $snapshotA[$oldRow]
    $ws.Cells.Item($newRow, 2).Value = $snapshotB[$oldRow]
    $ws.Cells.Item($newRow, 7).Value = $snapshotG[$oldRow]
}

# ---------------------------------------------------------------------------
# Step 3: rows 5 and 6 - swap the full record (every column, plus rewriting
# the HYPERLINK formulas so they keep pointing at "their own" Beteckning).
# ---------------------------------------------------------------------------
function Get-LinkFormula($column, $beteckning) {
    switch ($column) {
        "S" { return '=HYPERLINK("https://klasma.github.io/Logging_1481/artfynd/' + $beteckning + ' artfynd.xlsx", "' + $beteckning + '")' }
        "T" { return '=HYPERLINK("https://klasma.github.io/Logging_1481/kartor/' + $beteckning + ' karta.png", "' + $beteckning + '")' }
        "V" { return '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomål/' + $beteckning + ' FSC-klagomål.docx", "' + $beteckning + '")' }
        "W" { return '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomålsmail/' + $beteckning + ' FSC-klagomål mail.docx", "' + $beteckning + '")' }
        "X" { return '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsyn/' + $beteckning + ' tillsynsbegäran.docx", "' + $beteckning + '")' }
        "Y" { return '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsynsmail/' + $beteckning + ' tillsynsbegäran mail.docx", "' + $beteckning + '")' }
        "Z" { return '=HYPERLINK("https://klasma.github.io/Logging_1481/fåglar/' + $beteckning + ' prioriterade fågelarter.docx", "' + $beteckning + '")' }
    }
}

# Snapshot rows 5 and 6 completely (columns A, B, D, E, G..Q, R) before swapping.
$cols569 = @(1, 2, 4, 5, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18)
$row5 = @{}
$row6 = @{}
foreach ($c in $cols569) {
    $row5[$c] = $ws.Cells.Item(5, $c).Value2
    $row6[$c] = $ws.Cells.Item(6, $c).Value2
}
$beteckning5 = $ws.Cells.Item(5, 1).Value2
$beteckning6 = $ws.Cells.Item(6, 1).Value2

foreach ($c in $cols569) {
    $ws.Cells.Item(5, $c).Value = $row6[$c]
    $ws.Cells.Item(6, $c).Value = $row5[$c]
}

# Rewrite the HYPERLINK formulas for both rows against their new Beteckning.
$linkCols = @("S", "T", "V", "W", "X", "Y")
foreach ($col in $linkCols) {
    $ws.Range($col + "5").Formula = Get-LinkFormula $col $beteckning6
    $ws.Range($col + "6").Formula = Get-LinkFormula $col $beteckning5
}

# Row 5 (now "A 49789-2023") has no bird-survey link; row 6 (now "A 393-2025")
# gains the one that used to live on row 5.
$ws.Range("Z5").ClearContents()
$ws.Range("Z6").Formula = Get-LinkFormula "Z" $beteckning5
